# Auto-generated edit script: apply numeric cell updates per commit diff
# Source: Marilith_Profits.xlsx canonical OOXML diff

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 8 (Leve Item ID 4565)
$ws.Range("H8").Value = 76.2
$ws.Range("I8").Value = 84.75
$ws.Range("J8").Value = 42
$ws.Range("K8").Value = 254.25
$ws.Range("L8").Value = 126
$ws.Range("M8").Value = -115.25
$ws.Range("N8").Value = -404
# Row 32 (Leve Item ID 5484)
$ws.Range("H32").Value = 1044.2
$ws.Range("I32").Value = 779.4
$ws.Range("J32").Value = 1132.4667
$ws.Range("K32").Value = 779.4
$ws.Range("L32").Value = 1132.4667
$ws.Range("M32").Value = -453.4
$ws.Range("N32").Value = -1784.4667
# Row 38 (Leve Item ID 4599)
$ws.Range("H38").Value = 226.125
$ws.Range("I38").Value = 226.125
$ws.Range("K38").Value = 678.375
$ws.Range("M38").Value = -306.375
# Row 39 (Leve Item ID 4603)
$ws.Range("H39").Value = 504.73334
$ws.Range("I39").Value = 380.22223
$ws.Range("J39").Value = 691.5
$ws.Range("K39").Value = 1140.66669
$ws.Range("L39").Value = 2074.5
$ws.Range("M39").Value = -844.66669
$ws.Range("N39").Value = -2666.5
# Row 58 (Leve Item ID 4606)
$ws.Range("H58").Value = 1174.4
$ws.Range("I58").Value = 249.14285
$ws.Range("K58").Value = 747.4285500000001
$ws.Range("M58").Value = -597.4285500000001
# Row 82 (Leve Item ID 12623)
$ws.Range("H82").Value = 357.8889
$ws.Range("I82").Value = 357.8889
$ws.Range("K82").Value = 1073.6667
$ws.Range("M82").Value = -667.6667
# Row 85 (Leve Item ID 12623)
$ws.Range("H85").Value = 357.8889
$ws.Range("I85").Value = 357.8889
$ws.Range("K85").Value = 1073.6667
$ws.Range("M85").Value = 330.3333
# Row 116 (Leve Item ID 27778)
$ws.Range("H116").Value = 4671.625
$ws.Range("I116").Value = 3549
$ws.Range("K116").Value = 3549
$ws.Range("M116").Value = -107

$ws = $wb.Worksheets.Item("ARM")
# Row 2 (Leve Item ID 27713)
$ws.Range("H2").Value = 549.44446
$ws.Range("I2").Value = 493.125
$ws.Range("K2").Value = 493.125
$ws.Range("M2").Value = -380.125
# Row 97 (Leve Item ID 19941)
$ws.Range("H97").Value = 777.6923
$ws.Range("I97").Value = 758.25
$ws.Range("K97").Value = 758.25
$ws.Range("M97").Value = -262.25
# Row 116 (Leve Item ID 27713)
$ws.Range("H116").Value = 549.44446
$ws.Range("I116").Value = 493.125
$ws.Range("K116").Value = 493.125
$ws.Range("M116").Value = 1800.875
# Row 122 (Leve Item ID 36168)
$ws.Range("H122").Value = 3071.2856
$ws.Range("I122").Value = 1196
$ws.Range("K122").Value = 3588
$ws.Range("M122").Value = -1138

$ws = $wb.Worksheets.Item("BSM")
# Row 3 (Leve Item ID 27713)
$ws.Range("H3").Value = 549.44446
$ws.Range("I3").Value = 493.125
$ws.Range("K3").Value = 493.125
$ws.Range("M3").Value = -379.125
# Row 20 (Leve Item ID 14149)
$ws.Range("H20").Value = 1555.75
$ws.Range("I20").Value = 1546.5
$ws.Range("K20").Value = 1546.5
$ws.Range("M20").Value = -1299.5
# Row 94 (Leve Item ID 19939)
$ws.Range("H94").Value = 1083.5385
$ws.Range("I94").Value = 994.1818
$ws.Range("K94").Value = 994.1818
$ws.Range("M94").Value = -543.1818

$ws = $wb.Worksheets.Item("CRP")
# Row 16 (Leve Item ID 27691)
$ws.Range("H16").Value = 999.5
$ws.Range("I16").Value = 999.5
$ws.Range("K16").Value = 999.5
$ws.Range("M16").Value = -712.5
# Row 59 (Leve Item ID 1942)
$ws.Range("H59").Value = 29285.572
$ws.Range("I59").Value = 15000
$ws.Range("K59").Value = 15000
$ws.Range("M59").Value = -13855
# Row 107 (Leve Item ID 27689)
$ws.Range("H107").Value = 647
$ws.Range("I107").Value = 613.38464
$ws.Range("J107").Value = 719.8333
$ws.Range("K107").Value = 613.38464
$ws.Range("L107").Value = 719.8333
$ws.Range("M107").Value = 1306.61536
$ws.Range("N107").Value = -4559.8333
# Row 113 (Leve Item ID 27691)
$ws.Range("H113").Value = 999.5
$ws.Range("I113").Value = 999.5
$ws.Range("K113").Value = 999.5
$ws.Range("M113").Value = 1170.5

$ws = $wb.Worksheets.Item("CUL")
# Row 5 (Leve Item ID 43974)
$ws.Range("H5").Value = 827.6
$ws.Range("I5").Value = 827.6
$ws.Range("K5").Value = 2482.8
$ws.Range("M5").Value = -2370.8
# Row 99 (Leve Item ID 19817)
$ws.Range("H99").Value = 1900
$ws.Range("I99").Value = 1900
$ws.Range("K99").Value = 5700
$ws.Range("M99").Value = -3454
# Row 114 (Leve Item ID 27865)
$ws.Range("H114").Value = 0
$ws.Range("I114").Value = 0
$ws.Range("K114").Value = 0
$ws.Range("M114").ClearContents()
# Row 117 (Leve Item ID 27870)
$ws.Range("H117").Value = 174.5
$ws.Range("I117").Value = 174.5
$ws.Range("K117").Value = 523.5
$ws.Range("M117").Value = 2918.5
# Row 135 (Leve Item ID 43974)
$ws.Range("H135").Value = 827.6
$ws.Range("I135").Value = 827.6
$ws.Range("K135").Value = 7448.400000000001
$ws.Range("M135").Value = -4913.400000000001

$ws = $wb.Worksheets.Item("GSM")
# Row 39 (Leve Item ID 18264)
$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()
# Row 45 (Leve Item ID 27225)
$ws.Range("H45").Value = 20000
$ws.Range("J45").Value = 20000
$ws.Range("L45").Value = 20000
$ws.Range("N45").Value = -21118
# Row 102 (Leve Item ID 36169)
$ws.Range("H102").Value = 949.5
$ws.Range("I102").Value = 949.5
$ws.Range("K102").Value = 949.5
$ws.Range("M102").Value = 672.5

$ws = $wb.Worksheets.Item("LTW")
# Row 7 (Leve Item ID 36249)
$ws.Range("H7").Value = 13708.5625
$ws.Range("I7").Value = 12475.818
$ws.Range("J7").Value = 16420.6
$ws.Range("K7").Value = 12475.818
$ws.Range("L7").Value = 16420.6
$ws.Range("M7").Value = -12363.818
$ws.Range("N7").Value = -16644.6
# Row 98 (Leve Item ID 18379)
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()
# Row 107 (Leve Item ID 38752)
$ws.Range("H107").Value = 3997.5
$ws.Range("I107").Value = 3997.5
$ws.Range("K107").Value = 3997.5
$ws.Range("M107").Value = -2077.5
# Row 126 (Leve Item ID 36249)
$ws.Range("H126").Value = 13708.5625
$ws.Range("I126").Value = 12475.818
$ws.Range("J126").Value = 16420.6
$ws.Range("K126").Value = 37427.454
$ws.Range("L126").Value = 49261.8
$ws.Range("M126").Value = -34957.454
$ws.Range("N126").Value = -54201.8

$ws = $wb.Worksheets.Item("WVR")
# Row 81 (Leve Item ID 12596)
$ws.Range("H81").Value = 1134.8
$ws.Range("I81").Value = 1134.8
$ws.Range("K81").Value = 2269.6
$ws.Range("M81").Value = -1208.6
# Row 84 (Leve Item ID 12596)
$ws.Range("H84").Value = 1134.8
$ws.Range("I84").Value = 1134.8
$ws.Range("K84").Value = 11348
$ws.Range("M84").Value = -6044
# Row 107 (Leve Item ID 27746)
$ws.Range("H107").Value = 2098.0435
$ws.Range("I107").Value = 2347.6428
$ws.Range("J107").Value = 1709.7778
$ws.Range("K107").Value = 7042.928400000001
$ws.Range("L107").Value = 5129.3334
$ws.Range("M107").Value = -5122.928400000001
$ws.Range("N107").Value = -8969.3334
# Row 136 (Leve Item ID 44031)
$ws.Range("H136").Value = 3767.8333
$ws.Range("I136").Value = 4858
$ws.Range("J136").Value = 1587.5
$ws.Range("K136").Value = 14574
$ws.Range("L136").Value = 4762.5
$ws.Range("M136").Value = -12024
$ws.Range("N136").Value = -9862.5
